{"js": "// Remove the trailing boilerplate footer paragraphs that were dropped from\n// the Jekyll-generated page: the blank spacer paragraph, the\n// \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph, and the\n// \"\u00a9 2020 . Contact: ...\" paragraph. These three paragraphs immediately\n// follow the \"LOQ4038: ...\" requirements paragraph, which must stay intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that contains the LOQ4038 requirement text; the\n// three paragraphs to delete are the ones right after it.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"LOQ4038\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the LOQ4038 requirements paragraph.\");\n}\n\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (\n    text.trim() === \"\" ||\n    text.indexOf(\"Ver no Jupiter\") !== -1 ||\n    text.indexOf(\"Powered by Jekyll\") !== -1\n  ) {\n    toDelete.push(paragraphs.items[i]);\n    // Stop once we've collected the copyright/footer paragraph - everything\n    // that follows (trailing blank / page-break paragraphs) must remain.\n    if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n      break;\n    }\n  } else {\n    break;\n  }\n}\n\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing boilerplate footer paragraphs that were dropped from\n# the Jekyll-generated page: the blank spacer paragraph, the\n# \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph, and the\n# \"\u00a9 2020 . Contact: ...\" paragraph. These three paragraphs immediately\n# follow the \"LOQ4038: ...\" requirements paragraph, which must stay intact.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph containing the LOQ4038 requirement text.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($text -like \"*LOQ4038*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the LOQ4038 requirements paragraph.\"\n}\n\n# Collect the paragraphs right after the anchor that make up the footer:\n# a blank paragraph, the \"Ver no Jupiter ...\" paragraph, and the\n# \"... Powered by Jekyll ...\" copyright paragraph. Stop as soon as the\n# copyright paragraph has been collected so trailing paragraphs are kept.\n$toDelete = @()\nfor ($i = $anchorIndex + 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    $trimmed = $text.Trim()\n    if ($trimmed -eq \"\" -or $text -like \"*Ver no Jupiter*\" -or $text -like \"*Powered by Jekyll*\") {\n        $toDelete += $d.Paragraphs.Item($i)\n        if ($text -like \"*Powered by Jekyll*\") {\n            break\n        }\n    } else {\n        break\n    }\n}\n\n# Delete from the end backwards so earlier paragraph ranges stay valid.\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $toDelete[$j].Range.Delete()\n}\n"}
